# School transports' templates phrase update
#
# Adds a new "Body Text Indent" paragraph - "(μέσω της σχολικής μονάδας)" -
# right after the "Αναφερόμενους εκπαιδευτικούς" paragraph that closes the
# ΚΟΙΝΟΠΟΙΗΣΗ (distribution) list at the end of the document body.

$d = $word.ActiveDocument

# The target paragraph is the very last paragraph of the main body
# (immediately before the final section break).
$lastPara = $d.Paragraphs.Last
$tailRange = $lastPara.Range

# Create a brand new, empty paragraph right after it.
$tailRange.InsertParagraphAfter()

# That new paragraph is now the last paragraph in the document; grab its
# range so we can stamp it with the exact formatting/content from the diff.
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range

# Build the paragraph precisely: style "a8" (Body Text Indent), a left +
# first-line indent of 360 twips (0.25"/18pt each), and a single run in
# Calibri/MS Mincho (eastAsia) / Times New Roman (cs) at 11pt (sz/szCs 22
# half-points) carrying the new phrase.
$openXmlNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$newParaXml = "<w:p xmlns:w='$openXmlNs'>" +
    "<w:pPr>" +
        "<w:pStyle w:val='a8'/>" +
        "<w:ind w:left='360' w:firstLine='360'/>" +
    "</w:pPr>" +
    "<w:r>" +
        "<w:rPr>" +
            "<w:rFonts w:ascii='Calibri' w:eastAsia='MS Mincho;ＭＳ 明朝' w:hAnsi='Calibri' w:cs='Times New Roman'/>" +
            "<w:sz w:val='22'/>" +
            "<w:szCs w:val='22'/>" +
        "</w:rPr>" +
        "<w:t>(μέσω της σχολικής μονάδας)</w:t>" +
    "</w:r>" +
"</w:p>"

$newRange.InsertXML($newParaXml)
